$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 189
$ws.Range("I2").Value = 462
$ws.Range("J2").Value = 1989
$ws.Range("K2").Value = 14
$ws.Range("L2").Value = 535
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = 323
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 25
$ws.Range("S2").Value = 212
$ws.Range("T2").Value = 369
$ws.Range("U2").Value = 22
$ws.Range("V2").Value = 3061
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 3011
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 34
$ws.Range("AA2").Value = 17
